$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 98122.5
$ws.Range("J87").Value = 98122.5
$ws.Range("L87").Value = 98122.5
$ws.Range("N87").Value = -100618.5
$ws.Range("H90").Value = 98122.5
$ws.Range("J90").Value = 98122.5
$ws.Range("L90").Value = 294367.5
$ws.Range("N90").Value = -306847.5
$ws.Range("H125").Value = 3303.7896
$ws.Range("I125").Value = 3367.3845
$ws.Range("J125").Value = 3166
$ws.Range("K125").Value = 30306.4605
$ws.Range("L125").Value = 28494
$ws.Range("M125").Value = -27846.4605
$ws.Range("N125").Value = -33414
$ws.Range("H134").Value = 68194
$ws.Range("J134").Value = 68194
$ws.Range("L134").Value = 68194
$ws.Range("N134").Value = -78334
$ws.Range("H139").Value = 79780
$ws.Range("J139").Value = 79780
$ws.Range("L139").Value = 79780
$ws.Range("N139").Value = -90060
$ws.Range("H140").Value = 85390
$ws.Range("J140").Value = 85390
$ws.Range("L140").Value = 85390
$ws.Range("N140").Value = -95750

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4160.0244
$ws.Range("I132").Value = 4064.05
$ws.Range("K132").Value = 12192.15
$ws.Range("M132").Value = -9662.150000000001
$ws.Range("H134").Value = 94999.5
$ws.Range("J134").Value = 94999.5
$ws.Range("L134").Value = 94999.5
$ws.Range("N134").Value = -105139.5
$ws.Range("H139").Value = 75000
$ws.Range("J139").Value = 75000
$ws.Range("L139").Value = 75000
$ws.Range("N139").Value = -85280
$ws.Range("H141").Value = 52999.668
$ws.Range("J141").Value = 52999.668
$ws.Range("L141").Value = 52999.668
$ws.Range("N141").Value = -63359.668

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 89999
$ws.Range("J57").Value = 89999
$ws.Range("L57").Value = 89999
$ws.Range("N57").Value = -91439
$ws.Range("H134").Value = 3119.6924
$ws.Range("I134").Value = 2399.4285
$ws.Range("K134").Value = 7198.2855
$ws.Range("M134").Value = -4663.2855
$ws.Range("H135").Value = 74999
$ws.Range("J135").Value = 74999
$ws.Range("L135").Value = 74999
$ws.Range("N135").Value = -85139
$ws.Range("H136").Value = 89999
$ws.Range("J136").Value = 89999
$ws.Range("L136").Value = 89999
$ws.Range("N136").Value = -100199
$ws.Range("H137").Value = 39890
$ws.Range("J137").Value = 39890
$ws.Range("L137").Value = 39890
$ws.Range("N137").Value = -50090
$ws.Range("H138").Value = 89694.5
$ws.Range("J138").Value = 89694.5
$ws.Range("L138").Value = 89694.5
$ws.Range("N138").Value = -99974.5
$ws.Range("H139").Value = 71568.664
$ws.Range("J139").Value = 66998.5
$ws.Range("L139").Value = 66998.5
$ws.Range("N139").Value = -77278.5

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4122.853
$ws.Range("I31").Value = 3073.9473
$ws.Range("K31").Value = 3073.9473
$ws.Range("M31").Value = -2778.9473
$ws.Range("H34").Value = 4122.853
$ws.Range("I34").Value = 3073.9473
$ws.Range("K34").Value = 3073.9473
$ws.Range("M34").Value = -2871.9473
$ws.Range("H96").Value = 30964.6
$ws.Range("J96").Value = 30964.6
$ws.Range("L96").Value = 30964.6
$ws.Range("N96").Value = -36456.6
$ws.Range("H138").Value = 90333
$ws.Range("J138").Value = 90333
$ws.Range("L138").Value = 90333
$ws.Range("N138").Value = -100613
$ws.Range("H140").Value = 78330.664
$ws.Range("J140").Value = 92496
$ws.Range("L140").Value = 92496
$ws.Range("N140").Value = -102856

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 232.27272
$ws.Range("J17").Value = 805
$ws.Range("L17").Value = 2415
$ws.Range("N17").Value = -2753

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5694.88
$ws.Range("J122").Value = 7338.077
$ws.Range("L122").Value = 22014.231
$ws.Range("N122").Value = -26914.231
$ws.Range("H133").Value = 79885
$ws.Range("J133").Value = 79885
$ws.Range("L133").Value = 79885
$ws.Range("N133").Value = -90005
$ws.Range("H135").Value = 84685
$ws.Range("J135").Value = 84685
$ws.Range("L135").Value = 84685
$ws.Range("N135").Value = -94825
$ws.Range("H138").Value = 69997.5
$ws.Range("J138").Value = 69997.5
$ws.Range("L138").Value = 69997.5
$ws.Range("N138").Value = -80277.5
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 15388.625
$ws.Range("I7").Value = 12533
$ws.Range("J7").Value = 20148
$ws.Range("K7").Value = 12533
$ws.Range("L7").Value = 20148
$ws.Range("M7").Value = -12421
$ws.Range("N7").Value = -20372
$ws.Range("H126").Value = 15388.625
$ws.Range("I126").Value = 12533
$ws.Range("J126").Value = 20148
$ws.Range("K126").Value = 37599
$ws.Range("L126").Value = 60444
$ws.Range("M126").Value = -35129
$ws.Range("N126").Value = -65384
$ws.Range("H132").Value = 22718.2
$ws.Range("I132").Value = 30615.08
$ws.Range("K132").Value = 91845.24000000001
$ws.Range("M132").Value = -89315.24000000001
$ws.Range("H133").Value = 44999
$ws.Range("J133").Value = 44999
$ws.Range("L133").Value = 44999
$ws.Range("N133").Value = -50059
$ws.Range("H136").Value = 3911
$ws.Range("I136").Value = 3503.5
$ws.Range("K136").Value = 10510.5
$ws.Range("M136").Value = -7960.5
$ws.Range("H139").Value = 89999
$ws.Range("J139").Value = 89999
$ws.Range("L139").Value = 89999
$ws.Range("N139").Value = -100279
$ws.Range("H141").Value = 90000
$ws.Range("J141").Value = 90000
$ws.Range("L141").Value = 90000
$ws.Range("N141").Value = -100360

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I126").Value = 4748.85
$ws.Range("J126").Value = 5903.25
$ws.Range("K126").Value = 14246.55
$ws.Range("L126").Value = 17709.75
$ws.Range("M126").Value = -11776.55
$ws.Range("N126").Value = -22649.75
$ws.Range("H132").Value = 1309.6086
$ws.Range("I132").Value = 1204.4445
$ws.Range("J132").Value = 1688.2
$ws.Range("K132").Value = 3613.3335
$ws.Range("L132").Value = 5064.6
$ws.Range("M132").Value = -1083.3335
$ws.Range("N132").Value = -10124.6
$ws.Range("H137").Value = 63398.6
$ws.Range("J137").Value = 63398.6
$ws.Range("L137").Value = 63398.6
$ws.Range("N137").Value = -73598.60000000001
$ws.Range("H138").Value = 105999.5
$ws.Range("J138").Value = 105999.5
$ws.Range("L138").Value = 105999.5
$ws.Range("N138").Value = -116279.5
$ws.Range("H141").Value = 69990
$ws.Range("J141").Value = 69990
$ws.Range("L141").Value = 69990
$ws.Range("N141").Value = -80350
